$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared text so it reads "APPLE iPhone SE (White, 128 GB)"
$ws.Range("A1").Value2 = "APPLE iPhone SE (White, 128 GB)"

# The same text used to appear diagonally (B2, C3, D4, ... AE31).
# Move every one of those values into column A (A2..A31) and clear the old cell.
for ($i = 2; $i -le 31; $i++) {
    $oldCell = $ws.Cells.Item($i, $i)
    $newCell = $ws.Cells.Item($i, 1)
    $newCell.Value2 = $ws.Range("A1").Value2
    $oldCell.ClearContents()
}
